$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.57824772385901
$ws.Range("C2").Value = 0.125256517287666
$ws.Range("B3").Value = 0.0468565894924138
$ws.Range("C3").Value = 0.05928365352867

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.0046720337807
$ws.Range("C2").Value = 0.150216926418198
$ws.Range("B3").Value = -0.992824174200785
$ws.Range("C3").Value = 0.0678852275553906

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.99340653773269
$ws.Range("C2").Value = 0.0976327313995331
$ws.Range("B3").Value = 0.535253073568705
$ws.Range("C3").Value = 0.0774279575833975

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.3636120790641
$ws.Range("C2").Value = 0.126087808930259
$ws.Range("B3").Value = -0.0125577421591243
$ws.Range("C3").Value = 0.00856252069692879

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0156891951230354
$ws.Range("B2").Value = -0.00425682222495584
$ws.Range("A3").Value = -0.00425682222495584
$ws.Range("B3").Value = 0.00351455157570738

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0225651249825304
$ws.Range("B2").Value = -0.00802013418412126
$ws.Range("A3").Value = -0.00802013418412126
$ws.Range("B3").Value = 0.00460840412024716

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.00953215024053337
$ws.Range("B2").Value = 0.0017118540801892
$ws.Range("A3").Value = 0.0017118540801892
$ws.Range("B3").Value = 0.00599508861553641

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0158981355608334
$ws.Range("B2").Value = -0.000677055669639301
$ws.Range("A3").Value = -0.000677055669639301
$ws.Range("B3").Value = 0.0000733167606853339
